# Applies the LOM3114.xlsx content update:
#  - inserts a new row (row 13) to hold the teacher-name value on its own
#    (previously attached to "Objetivos:")
#  - fills in newly-authored long-form text for Objetivos / Programa resumido /
#    Programa / Bibliografia
# which together reproduce the row shift + text changes seen in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13; this pushes the old rows 13-23 down
# to 14-24 and leaves row 13 empty, matching the structural shift in the diff.
$ws.Rows(13).Insert()

# Row 10 ("Objetivos:") previously held the teacher's name in B10/C10; it now
# holds the long-form course-objectives paragraph in B10 only (C10 cleared).
$ws.Range("B10").Value = "Utilização de conceitos básicos da estatística para estudar influência de variáveis independentes sobre variáveis dependentes (respostas) em Processos da Engenharia . Utilização de ferramentas de planejamento experimental, dimensionamento das atividades experimentais de pesquisa e atividades industriais, minimizando tempo e custos, identificando quais são as variáveis de processo que podem influenciar significativamente nos sistemas estudados. Permitir a utilização de ferramentas estatísticas usando planilhas eletrônicas, para comparar metodologias e resultados em estudo de casos reais em Engenharia."
$ws.Range("C10").Value = ""

# New row 13 (blank label in A) now carries the teacher-name value that used
# to sit on the "Objetivos:" row.
$ws.Range("B13").Value = "5840521 - Rosa Ana Conte"
$ws.Range("C13").Value = "5840521 - Rosa Ana Conte"

# "Programa resumido:" (now row 14) gets its new long-form summary text,
# replacing the placeholder "Semestral".
$ws.Range("B14").Value = "Trabalho em planilhas eletrônicas e Estudo de Casos no Excel e Minitab; Introdução à Estatística Descritiva; Estatística de Inferência usando planilhas eletrônicas; distribuições amostrais; intervalos de confiança; testes de hipóteses ; testes ANOVA; estudo de casos em engenharia, meio ambiente, agricultura, gerenciamento de resíduos, dentre outros."
$ws.Range("C14").Value = "Trabalho em planilhas eletrônicas e Estudo de Casos no Excel e Minitab; Introdução à Estatística Descritiva; Estatística de Inferência usando planilhas eletrônicas; distribuições amostrais; intervalos de confiança; testes de hipóteses ; testes ANOVA; estudo de casos em engenharia, meio ambiente, agricultura, gerenciamento de resíduos, dentre outros."

# "Programa:" (now row 16) gets its new long-form syllabus text, replacing
# the placeholder date value.
$ws.Range("B16").Value = "O papel da estatística na Engenharia: métodos de coleta de dados. Trabalho em planilhas eletrônicas em Excel •Revisão de conceitos estatísticos fundamentais da estatística descritiva: população, amostra, tipos de erros associados a medidas experimentais.•Distribuições amostrais: distribuição normal, normal padronizada, de Student;•Estatística de inferência: estimativas, intervalos de confiança. •Testes de hipóteses para média aritmética e para duas populações: testes t uni- e bilateral, teste F, rejeição de valor suspeito; teste para independência ou homogeneidade da população • Análise de Variância (ANOVA): aplicações a problemas experimentais: fator único e 2 fatores; identificação de fatores significativos nos experimentos.•Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta Utilização de Minitab no planejamento e tomada de decisão de problemas experimentais. •Os conceitos desenvolvidos serão aplicados no estudo de casos reais nas áreas de engenharia, finanças, meio ambiente, agricultura, gerenciamento de resíduos, dentre outros."
$ws.Range("C16").Value = "O papel da estatística na Engenharia: métodos de coleta de dados. Trabalho em planilhas eletrônicas em Excel •Revisão de conceitos estatísticos fundamentais da estatística descritiva: população, amostra, tipos de erros associados a medidas experimentais.•Distribuições amostrais: distribuição normal, normal padronizada, de Student;•Estatística de inferência: estimativas, intervalos de confiança. •Testes de hipóteses para média aritmética e para duas populações: testes t uni- e bilateral, teste F, rejeição de valor suspeito; teste para independência ou homogeneidade da população • Análise de Variância (ANOVA): aplicações a problemas experimentais: fator único e 2 fatores; identificação de fatores significativos nos experimentos.•Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta Utilização de Minitab no planejamento e tomada de decisão de problemas experimentais. •Os conceitos desenvolvidos serão aplicados no estudo de casos reais nas áreas de engenharia, finanças, meio ambiente, agricultura, gerenciamento de resíduos, dentre outros."

# "Metodo:" (now row 19) keeps the "Estudo de Casos..." text that used to be
# one row up (on "Criterio:").
$ws.Range("B19").Value = "Estudo de Casos, Aulas expositivas e em laboratório computacional, trabalhos em grupo e exercícios comentados."
$ws.Range("C19").Value = "Estudo de Casos, Aulas expositivas e em laboratório computacional, trabalhos em grupo e exercícios comentados."

# "Criterio:" (now row 20) keeps the "Media aritmetica..." text.
$ws.Range("B20").Value = "Média aritmética de trabalhos propostos ao longo do curso (40%) e avaliação individual final (60%)."
$ws.Range("C20").Value = "Média aritmética de trabalhos propostos ao longo do curso (40%) e avaliação individual final (60%)."

# "Norma de recuperacao:" (now row 21) keeps the "Nao havera exame..." text.
$ws.Range("B21").Value = "Não haverá exame de recuperação."
$ws.Range("C21").Value = "Não haverá exame de recuperação."

# "Bibliografia:" (now row 22) gets the new long-form bibliography text,
# replacing the "Nao havera exame de recuperacao." placeholder.
$ws.Range("B22").Value = "•Levine, D.M et al. Estatística: teoria e Aplicações usando MicrosoftTM Excel em Português, 6a ed, Rio de Janeiro:LTC, 2012. •Mann, P.S. Introdução à Estatística, 8a ed, Rio deJaneiro:LTC, 2015. •Webster, A.L. Estatística Aplicada à Administração e Economia, São Paulo:McGraw Hill, 2007.•Johnson, R. e Kuby, P. ESTAT, São Paulo:Cengage Learning, 2014.•Barros Neto, B. , Scarminio, I.S. e Bruns, R.E. Planejamento e Otimização de Experimentos, 2a. ed, Campinas: Editora da UNICAMP, 1995.•Miller, JC and Miller, JN Statistical for Analytical Chemistry, Chichester: Ellishor Wood Ltd. 1988.•https://www.real-statistics.com• Kiernan, D. Natural Resources Biometrics: https://milnepublishing. genesco.edu/natural-resources-biometrics"
$ws.Range("C22").Value = "•Levine, D.M et al. Estatística: teoria e Aplicações usando MicrosoftTM Excel em Português, 6a ed, Rio de Janeiro:LTC, 2012. •Mann, P.S. Introdução à Estatística, 8a ed, Rio deJaneiro:LTC, 2015. •Webster, A.L. Estatística Aplicada à Administração e Economia, São Paulo:McGraw Hill, 2007.•Johnson, R. e Kuby, P. ESTAT, São Paulo:Cengage Learning, 2014.•Barros Neto, B. , Scarminio, I.S. e Bruns, R.E. Planejamento e Otimização de Experimentos, 2a. ed, Campinas: Editora da UNICAMP, 1995.•Miller, JC and Miller, JN Statistical for Analytical Chemistry, Chichester: Ellishor Wood Ltd. 1988.•https://www.real-statistics.com• Kiernan, D. Natural Resources Biometrics: https://milnepublishing. genesco.edu/natural-resources-biometrics"
